# Refresh the cryptocurrency snapshot table (columns D = Price,
# E = Volume(1h)) with the latest scraped figures.
#
# Source values are always plain text (see original inlineStr cells),
# so values that look like plain decimal numbers are written with a
# leading apostrophe to force Excel to keep them as text instead of
# silently converting them to numbers. The cell style is reset back
# to "Normal" afterwards so no stray text-format style is left behind.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.080.70"
$ws.Range("E2").Value = "  +1.34%  "
$ws.Range("D3").Value = "2.305.92"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "'310.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -1.78%  "
$ws.Range("D6").Value = "'106.09"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +2.11%  "
$ws.Range("E7").Value = "  -0.28%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("E9").Value = "  -0.05%  "
$ws.Range("D10").Value = "'39.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.83%  "
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "'8.31"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.12%  "
$ws.Range("D13").Value = "'0.108"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.20%  "
$ws.Range("E14").Value = "  -0.25%  "
$ws.Range("D15").Value = "'15.31"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.21%  "
$ws.Range("D16").Value = "2.656.07"
$ws.Range("E16").Value = "  +0.06%  "
$ws.Range("D17").Value = "2.303.26"
$ws.Range("E17").Value = "  +0.07%  "
$ws.Range("D18").Value = "42.741.12"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  -3.32%  "
$ws.Range("D21").Value = "'13.35"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.20%  "
$ws.Range("D22").Value = "'73.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.44%  "
$ws.Range("D23").Value = "'3.48"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.90%  "
$ws.Range("D24").Value = "'268.39"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.07%  "
$ws.Range("E25").Value = "  +0.10%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'7.49"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +11.67%  "
$ws.Range("D28").Value = "'10.99"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.63%  "
$ws.Range("E29").Value = "  -0.85%  "
$ws.Range("D30").Value = "'38.18"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.84%  "
$ws.Range("D31").Value = "'22.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.27%  "
$ws.Range("D32").Value = "'165.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.11%  "
$ws.Range("E33").Value = "  -1.79%  "
$ws.Range("D34").Value = "'2.83"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +5.80%  "
$ws.Range("E35").Value = "  -0.63%  "
$ws.Range("E36").Value = "  -1.06%  "
$ws.Range("E37").Value = "  +1.00%  "
$ws.Range("D38").Value = "'0.0352"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.70%  "
$ws.Range("E39").Value = "  +2.93%  "
$ws.Range("E40").Value = "  -2.12%  "
$ws.Range("D41").Value = "'108.13"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +9.34%  "
$ws.Range("E42").Value = "  -3.14%  "
$ws.Range("D43").Value = "'71.68"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +2.06%  "
$ws.Range("E44").Value = "  +1.16%  "
$ws.Range("E45").Value = "  -0.34%  "
$ws.Range("D46").Value = "'12.39"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.08%  "
$ws.Range("D47").Value = "1.721.12"
$ws.Range("E47").Value = "  +5.34%  "
$ws.Range("D48").Value = "'111.56"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.73%  "
$ws.Range("D49").Value = "'76.90"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -4.87%  "
$ws.Range("E50").Value = "  +0.45%  "
$ws.Range("D51").Value = "'5.20"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.92%  "
